$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: AC1 "Wins", AD1 "Losses", AE1 "Ties" ---
# Copy the existing header formatting (from A1, style index 1: bold/bordered/centered)
# onto the three new header cells before writing their text.
$ws.Range("A1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# --- Data rows 2-41: team record (Wins=83, Losses=79, Ties=0) on every player row ---
$lastRow = 41
$ws.Range("AC2:AC$lastRow").Value = 83
$ws.Range("AD2:AD$lastRow").Value = 79
$ws.Range("AE2:AE$lastRow").Value = 0
